$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.495.67"
$ws.Range("E2").Value = "  -0.41%  "

$ws.Range("D3").Value = "3.108.31"
$ws.Range("E3").Value = "  +2.93%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "385.17"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.49"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.32%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.541"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.02%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.35"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.137"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0855"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.62%  "

$ws.Range("D13").Value = "3.592.52"
$ws.Range("E13").Value = "  +3.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.65"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.83"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "

$ws.Range("D16").Value = "3.107.89"
$ws.Range("E16").Value = "  +3.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.994"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.98"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.09%  "

$ws.Range("D19").Value = "51.539.62"
$ws.Range("E19").Value = "  -0.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.27"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.39"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.54%  "

$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.05"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "266.24"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.09"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.11"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.09"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.46%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.25"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.169"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "

$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.107"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.33"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "35.61"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.81%  "

$ws.Range("E34").Value = "  +4.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.07"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.32"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.85%  "

$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.37"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.296"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.85%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "129.19"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.39%  "

$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.80"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.47"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.67"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.47%  "

$ws.Range("E47").Value = "  +4.15%  "

$ws.Range("E48").Value = "  +1.13%  "

$ws.Range("D49").Value = "2.062.18"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").Value = "3.416.64"
$ws.Range("E50").Value = "  +3.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.913"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +16.90%  "
